# Updates ASIN report rows: refresh existing product metrics (rows 2-12)
# and append newly tracked Majestic Pure SKUs (rows 13-20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ r=2; A="Majestic Pure"; B="NG-IWXD-5PE7"; C="Majestic Pure Bulgarian Lavender Essential Oil, 100% Pure and Natural with Therapeutic Grade, Premium Quality Bulgarian Lavender Oil, 1 fl. oz."; D="Active"; E="B01FZRK3WW"; F="354"; G="25"; H="4.5"; I="63,169"; J="366"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=3; A="Majestic Pure"; B="Acnea3ScarWashCleanser"; C="MAJESTIC PURE Acnea3 Scar Acne Wash Foaming Facial Cleanser - Soothes Blemishes and Clears Pores, Made with Pure Peppermint and Tea Tree Essential Oils - for All Skin Types, 4 fl oz"; D="Active"; E="B07VXK2RX4"; F="210"; G=""; H="4.1"; I="153,219"; J="1,592"; K=""; L="Available"; M="Lose Q & A, Sub.Cat2" }
    @{ r=4; A="Majestic Pure"; B="6L-S5IM-ZW9C"; C="Majestic Pure Fractionated Coconut Oil, For Aromatherapy Relaxing Massage, Carrier Oil for Diluting Essential Oils, Hair & Skin Care Benefits, Moisturizer & Softener - 16 Ounces (Packaging May Vary)"; D="Active"; E="B00PMR3QF2"; F="19164"; G="394"; H="4.7"; I="583"; J="14"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=5; A="Majestic Pure"; B="NewPeppermintOil"; C="Majestic Pure Peppermint Essential Oil, Pure and Natural, Therapeutic Grade Peppermint Oil, 4 fl. oz."; D="Active"; E="B00PV15BPW"; F="12842"; G="193"; H="4.7"; I="3,519"; J="14"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=6; A="Majestic Pure"; B="BX-N0NF-ETC9"; C="MAJESTIC PURE Indian Healing Clay Powder, Deep Pore Cleansing Facial, Body and Hair Mask, Natural Sodium Bentonite Clay, 16oz"; D="Active"; E="B00Q96XGUU"; F="1785"; G="109"; H="4.3"; I="122,474"; J="2,391"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=7; A="Majestic Pure"; B="NewLemonOil"; C="Majestic Pure Lemon Oil, Therapeutic Grade, Premium Quality Lemon Oil, 4 Ounces"; D="Active"; E="B00QR6SS6O"; F="6769"; G="109"; H="4.6"; I="5,597"; J="26"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=8; A="Majestic Pure"; B="AJ-CD9K-SAOC"; C="Majestic Pure Rosehip Oil for Face, Nails, Hair and Skin, Pure & Natural, Cold Pressed Premium Rose Hip Seed Oil, 4 oz"; D="Active"; E="B00QR7FTLU"; F="4599"; G="121"; H="4.6"; I="14,832"; J="145"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=9; A="Majestic Pure"; B="UT-X84W-ZTPT"; C="Majestic Pure Moroccan Argan Oil for Hair, Face, Nails, Beard & Cuticles - for Men and Women - Pure & Natural, 4 fl. oz."; D="Active"; E="B00QVR0O6Q"; F="3580"; G="103"; H="4.6"; I="91,531"; J="868"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=10; A="Majestic Pure"; B="JB-RO97-9L06"; C="Majestic Pure Jojoba Oil for Hair and Skin, 4 fl. oz."; D="Active"; E="B00STVN68K"; F="859"; G="38"; H="4.7"; I="53,541"; J="504"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=11; A="Majestic Pure"; B="MY-ETEY-F5KK"; C="Majestic Pure Lavender Oil, Natural, Therapeutic Grade, Premium Quality Blend of Lavender Essential Oil, 4 fl. Oz"; D="Active"; E="B00TSTZQEY"; F="15736"; G="217"; H="4.6"; I="4,373"; J="19"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=12; A="Majestic Pure"; B="9D-BZC2-93AO"; C="MAJESTIC PURE Dead Sea Mud Mask - Natural Face and Skin Care for Women and Men - Best Black Facial Cleansing Clay for Blackhead, Whitehead, Acne and Pores - 8.8 fl. Oz"; D="Active"; E="B00UREAGU8"; F="5614"; G="126"; H="4.5"; I="4,324"; J="59"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=13; A="Majestic Pure"; B="LV-XBR2-FJDV"; C="Majestic Pure Castor Oil, 100% Natural Wonder Oil with Numerous Hair, Scalp, Skin and Nails Benefits - Packaging May Vary- 16 fl oz"; D="Active"; E="B00XE58NJ8"; F="838"; G="59"; H="4.5"; I="44,194"; J="420"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=14; A="Majestic Pure"; B="KO-M9OF-WEJA"; C="Majestic Pure Hair Loss Shampoo, Offers Natural Ingredient Based Effective Solution, Add Volume and Strengthen Hair, Sulfate Free, 14 DHT Blockers, for Men & Women - 16 fl Oz"; D="Active"; E="B016RQ8PRU"; F="1654"; G="90"; H="4"; I="27,836"; J="40"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=15; A="Majestic Pure"; B="48-62AK-XNDW"; C="Majestic Pure Black Pepper Essential Oil, Pure and Natural with Therapeutic Grade, Premium Quality Black Pepper Oil, 4 fl. oz."; D="Active"; E="B0172AL2PM"; F="178"; G="3"; H="4.5"; I="106,561"; J="652"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=16; A="Majestic Pure"; B="7W-UX75-3750"; C="Majestic Pure Clary Sage Oil, Premium Quality, 4 fl. oz."; D="Active"; E="B01767OTVK"; F="667"; G="25"; H="4.6"; I="58,813"; J="341"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=17; A="Majestic Pure"; B="G9-CLOG-DZRL"; C="Majestic Pure Juniper Oil, Premium Quality, 4 fl. oz"; D="Active"; E="B0176UQBJU"; F="473"; G="17"; H="4.3"; I=""; J=""; K=""; L="Available"; M="Lose Category, Sub. Cat, Sub.Cat2" }
    @{ r=18; A="Majestic Pure"; B="KO-I58S-WEWY"; C="Majestic Pure Myrrh Oil, Premium Quality, 4 fl Oz"; D="Active"; E="B0176YOIEQ"; F="1019"; G="23"; H="4.4"; I="41,122"; J="228"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=19; A="Majestic Pure"; B="2H-DX31-4SOC"; C="Majestic Pure Coconut Milk Body Scrub, Anti Cellulite & Exfoliator, Natural Skin Care Formula Helps with Stretch Marks, Eczema, Acne and Varicose Veins, 12 Oz"; D="Active"; E="B01BHQ69S2"; F="810"; G="18"; H="4.1"; I="29,576"; J="110"; K=""; L="Available"; M="Lose Sub.Cat2" }
    @{ r=20; A="Majestic Pure"; B="4J-YNLD-KLGK"; C="Majestic Pure Grapefruit Essential Oil, Pure and Natural, Therapeutic Grade Grapefruit Oil, 4 fl. oz."; D="Active"; E="B01BKALLBU"; F="945"; G="22"; H="4.6"; I=""; J=""; K=""; L="Available"; M="Lose Category, Sub. Cat, Sub.Cat2" }
)

foreach ($row in $rows) {
    $r = $row.r

    # A-E: plain text labels / codes / ASINs.
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E

    # F: Customer Reviews (numeric)
    $ws.Cells.Item($r, 6).Value = [double]$row.F

    # G: Q & A (numeric, blank on some rows)
    if ($row.G -eq "") {
        $ws.Cells.Item($r, 7).Value = ""
    } else {
        $ws.Cells.Item($r, 7).Value = [double]$row.G
    }

    # H: Reviews Rating (numeric)
    $ws.Cells.Item($r, 8).Value = [double]$row.H

    # I, J, K: Category / Sub. Cat / Sub.Cat2 are stored as literal text
    # even when they look like numbers (e.g. "63,169"), so force text
    # format before assigning to avoid Excel auto-converting them.
    $ws.Cells.Item($r, 9).NumberFormat = "@"
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).NumberFormat = "@"
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).NumberFormat = "@"
    $ws.Cells.Item($r, 11).Value = $row.K

    # L, M: plain text.
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
}
